$wb = $excel.ActiveWorkbook

# Data update for 2024-12-28 (adds new crime incident records reflected
# across the citywide totals, neighborhood summary, and per-neighborhood
# breakdown sheets).

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("K2").Value = 150
$ws.Range("C3").Value = 81
$ws.Range("C9").Value = 507
$ws.Range("D9").Value = 440
$ws.Range("E9").Value = 509
$ws.Range("F9").Value = 575
$ws.Range("G9").Value = 448
$ws.Range("H9").Value = 479
$ws.Range("I9").Value = 518
$ws.Range("J9").Value = 435
$ws.Range("B10").Value = 1410
$ws.Range("C10").Value = 1667
$ws.Range("D10").Value = 1888
$ws.Range("E10").Value = 2318
$ws.Range("F10").Value = 2208
$ws.Range("G10").Value = 922
$ws.Range("H10").Value = 636
$ws.Range("K10").Value = 709
$ws.Range("B11").Value = 1944
$ws.Range("C11").Value = 2339
$ws.Range("D11").Value = 2575
$ws.Range("E11").Value = 3067
$ws.Range("F11").Value = 3038
$ws.Range("G11").Value = 1619
$ws.Range("H11").Value = 1409
$ws.Range("I11").Value = 1748
$ws.Range("J11").Value = 1598
$ws.Range("K11").Value = 1661

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("D7").Value = 33
$ws.Range("J7").Value = 33
$ws.Range("E8").Value = 107
$ws.Range("D9").Value = 98
$ws.Range("E9").Value = 175
$ws.Range("J9").Value = 83

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("K2").Value = 8
$ws.Range("K10").Value = 57

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("C8").Value = 37
$ws.Range("C10").Value = 106

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("G8").Value = 57
$ws.Range("B9").Value = 231
$ws.Range("C9").Value = 350
$ws.Range("D9").Value = 537
$ws.Range("E9").Value = 698
$ws.Range("F9").Value = 562
$ws.Range("K9").Value = 113
$ws.Range("B10").Value = 280
$ws.Range("C10").Value = 411
$ws.Range("D10").Value = 615
$ws.Range("E10").Value = 788
$ws.Range("F10").Value = 649
$ws.Range("G10").Value = 255
$ws.Range("K10").Value = 207

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("G5").Value = 15
$ws.Range("G7").Value = 25

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("F7").Value = 57
$ws.Range("H7").Value = 8
$ws.Range("F8").Value = 99
$ws.Range("H8").Value = 25

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("F8").Value = 165
$ws.Range("H8").Value = 104
$ws.Range("D11").Value = 7
$ws.Range("K19").Value = 57
$ws.Range("D27").Value = 30
$ws.Range("D28").Value = 106
$ws.Range("F28").Value = 134
$ws.Range("G28").Value = 92
$ws.Range("D29").Value = 29
$ws.Range("D32").Value = 98
$ws.Range("E32").Value = 175
$ws.Range("J32").Value = 83
$ws.Range("D35").Value = 25
$ws.Range("C36").Value = 106
$ws.Range("C38").Value = 8
$ws.Range("I41").Value = 20
$ws.Range("G47").Value = 52
$ws.Range("G50").Value = 25
$ws.Range("J52").Value = 31
$ws.Range("B53").Value = 280
$ws.Range("C53").Value = 411
$ws.Range("D53").Value = 615
$ws.Range("E53").Value = 788
$ws.Range("F53").Value = 649
$ws.Range("G53").Value = 255
$ws.Range("K53").Value = 207
$ws.Range("E61").Value = 68
$ws.Range("H61").Value = 28
$ws.Range("J61").Value = 5
$ws.Range("F62").Value = 28
$ws.Range("F63").Value = 13
$ws.Range("F65").Value = 99
$ws.Range("H65").Value = 25
$ws.Range("E74").Value = 80
$ws.Range("F74").Value = 91
$ws.Range("B76").Value = 59
$ws.Range("E76").Value = 101
$ws.Range("F76").Value = 69
$ws.Range("F77").Value = 77
$ws.Range("F88").Value = 10
$ws.Range("F89").Value = 24
$ws.Range("I92").Value = 40
$ws.Range("C95").Value = 34
$ws.Range("D95").Value = 58
$ws.Range("B99").Value = 1944
$ws.Range("C99").Value = 2339
$ws.Range("D99").Value = 2575
$ws.Range("E99").Value = 3067
$ws.Range("F99").Value = 3038
$ws.Range("G99").Value = 1619
$ws.Range("H99").Value = 1409
$ws.Range("I99").Value = 1748
$ws.Range("J99").Value = 1598
$ws.Range("K99").Value = 1661

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("F6").Value = 17
$ws.Range("F7").Value = 24

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("I6").Value = 7
$ws.Range("I8").Value = 20

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("D7").Value = 25
$ws.Range("F8").Value = 73
$ws.Range("G8").Value = 41
$ws.Range("D9").Value = 106
$ws.Range("F9").Value = 134
$ws.Range("G9").Value = 92

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("G7").Value = 33
$ws.Range("G8").Value = 52

$ws = $wb.Worksheets.Item('Fuller Park')
$ws.Range("D7").Value = 5
$ws.Range("D9").Value = 29

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("E8").Value = 10
$ws.Range("B9").Value = 56
$ws.Range("F9").Value = 49
$ws.Range("B10").Value = 59
$ws.Range("E10").Value = 101
$ws.Range("F10").Value = 69

$ws = $wb.Worksheets.Item('River North')
$ws.Range("F5").Value = 10
$ws.Range("E6").Value = 73
$ws.Range("E7").Value = 80
$ws.Range("F7").Value = 91

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("I7").Value = 7
$ws.Range("I9").Value = 40

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("D6").Value = 15
$ws.Range("D7").Value = 25

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("J6").Value = 8
$ws.Range("J8").Value = 31

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range("D7").Value = 21
$ws.Range("D8").Value = 30

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range("F7").Value = 23
$ws.Range("F8").Value = 28

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("F9").Value = 54
$ws.Range("F10").Value = 77

$ws = $wb.Worksheets.Item('West Town')
$ws.Range("C6").Value = 28
$ws.Range("D6").Value = 46
$ws.Range("C7").Value = 34
$ws.Range("D7").Value = 58

$ws = $wb.Worksheets.Item('New City')
$ws.Range("F5").Value = 8
$ws.Range("F6").Value = 13

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("D6").Value = 6
$ws.Range("D7").Value = 7

$ws = $wb.Worksheets.Item('Greektown')
$ws.Range("C3").Value = 1
$ws.Range("C7").Value = 8

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("H7").Value = 50
$ws.Range("F8").Value = 112
$ws.Range("F9").Value = 165
$ws.Range("H9").Value = 104

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("F6").Value = 6
$ws.Range("F7").Value = 10
